$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 250; this shifts the existing rows 250-259 down to 251-260,
# matching the diff (a new weekly price record is inserted before the former row 250).
$ws.Rows.Item(250).Insert()

# Populate the newly inserted row 250 with the new record's data.
$ws.Cells.Item(250, 1).Value = 10
$ws.Cells.Item(250, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(250, 3).Value = "La Araucanía"
$ws.Cells.Item(250, 4).Value = 44753
$ws.Cells.Item(250, 5).Value = 9
$ws.Cells.Item(250, 6).Value = 100112039
$ws.Cells.Item(250, 7).Value = "Ciboulette"
$ws.Cells.Item(250, 8).Value = "Sin especificar"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 90
$ws.Cells.Item(250, 11).Value = 6000
$ws.Cells.Item(250, 12).Value = 7000
$ws.Cells.Item(250, 13).Value = 6444
$ws.Cells.Item(250, 14).Value = "$/docena de atados"
$ws.Cells.Item(250, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(250, 16).Value = 2148
$ws.Cells.Item(250, 17).Value = 3
$ws.Cells.Item(250, 18).Value = "Hortaliza"

# Make sure the new date cell uses the same custom date number format as the other cells in column D.
$ws.Cells.Item(250, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
